$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: add P1 and Q1, copying the header style (s=1) from an existing header cell ---
$ws.Range("O1").Copy()
$ws.Range("P1:Q1").PasteSpecial(-4122)
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15

# --- Data rows 2-25 ---
# Row 2
$ws.Range("B2").Value = 24.52631444203648
$ws.Range("C2").Value = 19.04901804811669
$ws.Range("D2").Value = 4.288684017742457
$ws.Range("E2").Value = 29.2869706272671
$ws.Range("F2").Value = 18.46192158008147
$ws.Range("H2").Value = 3.205158693868897
$ws.Range("I2").Value = 3.20525599313035
$ws.Range("J2").ClearContents()
$ws.Range("O2").Value = 0
$ws.Range("P2").Value = 0
$ws.Range("Q2").Value = 14.41630593039575

# Row 3
$ws.Range("B3").Value = 22.91929028977459
$ws.Range("C3").Value = 17.94061183013837
$ws.Range("D3").Value = 4.15026255036736
$ws.Range("E3").Value = 27.39456778111133
$ws.Range("F3").Value = 17.86194650852501
$ws.Range("H3").Value = 2.946021979672167
$ws.Range("I3").Value = 3.009413731153047
$ws.Range("J3").ClearContents()
$ws.Range("O3").Value = 0
$ws.Range("P3").Value = 0
$ws.Range("Q3").Value = 14.17742455656328

# Row 4
$ws.Range("B4").Value = 21.87328337603005
$ws.Range("C4").Value = 17.23026101314926
$ws.Range("D4").Value = 4.062780859921553
$ws.Range("E4").Value = 26.16811167537038
$ws.Range("F4").Value = 17.4957900591206
$ws.Range("H4").Value = 2.78160668458218
$ws.Range("I4").Value = 2.886299199571967
$ws.Range("J4").ClearContents()
$ws.Range("O4").Value = 0
$ws.Range("P4").Value = 0
$ws.Range("Q4").Value = 14.03809575590276

# Row 5
$ws.Range("B5").Value = 21.43188166716974
$ws.Range("C5").Value = 16.94920400894058
$ws.Range("D5").Value = 4.028706384070023
$ws.Range("E5").Value = 25.65183158974289
$ws.Range("F5").Value = 17.33325217044634
$ws.Range("H5").Value = 2.713101249232736
$ws.Range("I5").Value = 2.835859332754361
$ws.Range("J5").ClearContents()
$ws.Range("O5").Value = 0
$ws.Range("P5").Value = 0
$ws.Range("Q5").Value = 13.97213157276983

# Row 6
$ws.Range("B6").Value = 21.35745317030847
$ws.Range("C6").Value = 16.92142796309843
$ws.Range("D6").Value = 4.025682037953141
$ws.Range("E6").Value = 25.56483416297804
$ws.Range("F6").Value = 17.28911852458194
$ws.Range("H6").Value = 2.701467989658521
$ws.Range("I6").Value = 2.828207529254503
$ws.Range("J6").ClearContents()
$ws.Range("O6").Value = 0
$ws.Range("P6").Value = 0
$ws.Range("Q6").Value = 13.94783521824851

# Row 7
$ws.Range("B7").Value = 21.86678151261478
$ws.Range("C7").Value = 17.2782544451979
$ws.Range("D7").Value = 4.069565799390239
$ws.Range("E7").Value = 26.16045078440463
$ws.Range("F7").Value = 17.44715002475733
$ws.Range("H7").Value = 2.780220371406649
$ws.Range("I7").Value = 2.887359795471857
$ws.Range("J7").ClearContents()
$ws.Range("O7").Value = 0
$ws.Range("P7").Value = 0
$ws.Range("Q7").Value = 14.0006616968689

# Row 8
$ws.Range("B8").Value = 23.98374922097493
$ws.Range("C8").Value = 18.73708492387128
$ws.Range("D8").Value = 4.250852173765058
$ws.Range("E8").Value = 28.64691035016888
$ws.Range("F8").Value = 18.19570782938383
$ws.Range("H8").Value = 3.116316921085846
$ws.Range("I8").Value = 3.140047652784549
$ws.Range("J8").ClearContents()
$ws.Range("O8").Value = 0
$ws.Range("P8").Value = 0
$ws.Range("Q8").Value = 14.28511541946979

# Row 9
$ws.Range("B9").Value = 27.66824537335466
$ws.Range("C9").Value = 21.27434326780041
$ws.Range("D9").Value = 4.573118291163154
$ws.Range("E9").Value = 33.01452523321491
$ws.Range("F9").Value = 19.75299864164758
$ws.Range("H9").Value = 3.735367590512782
$ws.Range("I9").Value = 3.611985141512975
$ws.Range("J9").ClearContents()
$ws.Range("O9").Value = 0
$ws.Range("P9").Value = 0
$ws.Range("Q9").Value = 14.9659930966735

# Row 10
$ws.Range("B10").Value = 30.08969881375781
$ws.Range("C10").Value = 23.00217623967392
$ws.Range("D10").Value = 4.835908615770593
$ws.Range("E10").Value = 35.03852183181854
$ws.Range("F10").Value = 20.68522245521129
$ws.Range("H10").Value = 4.129536869478232
$ws.Range("I10").Value = 3.939775442244501
$ws.Range("J10").ClearContents()
$ws.Range("O10").Value = 0
$ws.Range("P10").Value = 0
$ws.Range("Q10").Value = 15.33064473205666

# Row 11
$ws.Range("B11").Value = 31.09754528126896
$ws.Range("C11").Value = 23.7509481832237
$ws.Range("D11").Value = 5.285367748218243
$ws.Range("E11").Value = 28.44207018043513
$ws.Range("F11").Value = 19.41062700068129
$ws.Range("H11").Value = 4.492686363854161
$ws.Range("I11").Value = 4.026397823387395
$ws.Range("J11").ClearContents()
$ws.Range("O11").Value = 0
$ws.Range("P11").Value = 0
$ws.Range("Q11").Value = 14.10346184430302

# Row 12
$ws.Range("B12").Value = 31.46003994495555
$ws.Range("C12").Value = 23.97306138816726
$ws.Range("D12").Value = 5.611338372245933
$ws.Range("E12").Value = 22.49534382238961
$ws.Range("F12").Value = 18.18018490245333
$ws.Range("H12").Value = 5.301114636066927
$ws.Range("I12").Value = 4.036352383828929
$ws.Range("J12").ClearContents()
$ws.Range("O12").Value = 0
$ws.Range("P12").Value = 0
$ws.Range("Q12").Value = 13.03009914932889

# Row 13
$ws.Range("B13").Value = 31.36001126536393
$ws.Range("C13").Value = 23.87246084525962
$ws.Range("D13").Value = 5.875420566203079
$ws.Range("E13").Value = 16.57826779130319
$ws.Range("F13").Value = 16.81684718730127
$ws.Range("H13").Value = 6.321874535282032
$ws.Range("I13").Value = 3.990327379479631
$ws.Range("J13").ClearContents()
$ws.Range("O13").Value = 0
$ws.Range("P13").Value = 0
$ws.Range("Q13").Value = 11.94702447354848

# Row 14
$ws.Range("B14").Value = 31.07990833678944
$ws.Range("C14").Value = 23.66791434137551
$ws.Range("D14").Value = 6.037018980021518
$ws.Range("E14").Value = 12.51407151035088
$ws.Range("F14").Value = 15.77948081056859
$ws.Range("H14").Value = 7.120997004698546
$ws.Range("I14").Value = 3.934376082290615
$ws.Range("J14").ClearContents()
$ws.Range("O14").Value = 0
$ws.Range("P14").Value = 0
$ws.Range("Q14").Value = 11.17789104184732

# Row 15
$ws.Range("B15").Value = 30.91110089188882
$ws.Range("C15").Value = 23.56401353955093
$ws.Range("D15").Value = 6.065376218933928
$ws.Range("E15").Value = 11.53784735484873
$ws.Range("F15").Value = 15.48355151738994
$ws.Range("H15").Value = 7.304739452644291
$ws.Range("I15").Value = 3.908720832596522
$ws.Range("J15").ClearContents()
$ws.Range("O15").Value = 0
$ws.Range("P15").Value = 0
$ws.Range("Q15").Value = 10.97941070898476

# Row 16
$ws.Range("B16").Value = 29.93972865017516
$ws.Range("C16").Value = 22.89778659810451
$ws.Range("D16").Value = 5.91695598565147
$ws.Range("E16").Value = 11.40079220748174
$ws.Range("F16").Value = 15.26308924785205
$ws.Range("H16").Value = 7.023652367670555
$ws.Range("I16").Value = 3.780498500792048
$ws.Range("J16").ClearContents()
$ws.Range("O16").Value = 0
$ws.Range("P16").Value = 0
$ws.Range("Q16").Value = 11.00000258495085

# Row 17
$ws.Range("B17").Value = 29.33335524780634
$ws.Range("C17").Value = 22.49134673905846
$ws.Range("D17").Value = 5.712032904846049
$ws.Range("E17").Value = 13.49393089274596
$ws.Range("F17").Value = 15.66491031130633
$ws.Range("H17").Value = 6.331781949282495
$ws.Range("I17").Value = 3.711235346739993
$ws.Range("J17").ClearContents()
$ws.Range("O17").Value = 0
$ws.Range("P17").Value = 0
$ws.Range("Q17").Value = 11.42603845118459

# Row 18
$ws.Range("B18").Value = 28.98815594547917
$ws.Range("C18").Value = 22.23687993014695
$ws.Range("D18").Value = 5.433055525278049
$ws.Range("E18").Value = 18.04966754262918
$ws.Range("F18").Value = 16.67421346875417
$ws.Range("H18").Value = 5.295175949569644
$ws.Range("I18").Value = 3.686294538833137
$ws.Range("J18").ClearContents()
$ws.Range("O18").Value = 0
$ws.Range("P18").Value = 0
$ws.Range("Q18").Value = 12.27652785648718

# Row 19
$ws.Range("B19").Value = 28.88555310593026
$ws.Range("C19").Value = 22.21234335816534
$ws.Range("D19").Value = 5.137946159416773
$ws.Range("E19").Value = 24.32584718494821
$ws.Range("F19").Value = 18.01095879176808
$ws.Range("H19").Value = 4.315711853095834
$ws.Range("I19").Value = 3.710300522297746
$ws.Range("J19").ClearContents()
$ws.Range("O19").Value = 0
$ws.Range("P19").Value = 0
$ws.Range("Q19").Value = 13.35529296422902

# Row 20
$ws.Range("B20").Value = 29.47103432876578
$ws.Range("C20").Value = 22.68727170978404
$ws.Range("D20").Value = 4.790412143687933
$ws.Range("E20").Value = 34.47749105733784
$ws.Range("F20").Value = 20.30784993647318
$ws.Range("H20").Value = 4.023774908018307
$ws.Range("I20").Value = 3.859958237290177
$ws.Range("J20").ClearContents()
$ws.Range("O20").Value = 0
$ws.Range("P20").Value = 0
$ws.Range("Q20").Value = 15.12155973123076

# Row 21
$ws.Range("B21").Value = 31.24823133926937
$ws.Range("C21").Value = 23.94496119903691
$ws.Range("D21").Value = 4.926019927894772
$ws.Range("E21").Value = 37.30559408676149
$ws.Range("F21").Value = 21.34808689034557
$ws.Range("H21").Value = 4.378860539700025
$ws.Range("I21").Value = 4.119120617262744
$ws.Range("J21").ClearContents()
$ws.Range("O21").Value = 0
$ws.Range("P21").Value = 0
$ws.Range("Q21").Value = 15.68174502682728

# Row 22
$ws.Range("B22").Value = 32.35835999905605
$ws.Range("C22").Value = 24.68646706035022
$ws.Range("D22").Value = 5.022453335923683
$ws.Range("E22").Value = 38.645784568976
$ws.Range("F22").Value = 21.97536135577538
$ws.Range("H22").Value = 4.5880490719886
$ws.Range("I22").Value = 4.280603779598847
$ws.Range("J22").ClearContents()
$ws.Range("O22").Value = 0
$ws.Range("P22").Value = 0
$ws.Range("Q22").Value = 16.01548221594701

# Row 23
$ws.Range("B23").Value = 31.7709066848655
$ws.Range("C23").Value = 24.24709859273765
$ws.Range("D23").Value = 4.963622951361446
$ws.Range("E23").Value = 37.93605119404339
$ws.Range("F23").Value = 21.68523465542879
$ws.Range("H23").Value = 4.477121239092268
$ws.Range("I23").Value = 4.192357459237603
$ws.Range("J23").ClearContents()
$ws.Range("O23").Value = 0
$ws.Range("P23").Value = 0
$ws.Range("Q23").Value = 15.87452321876289

# Row 24
$ws.Range("B24").Value = 29.44524700556518
$ws.Range("C24").Value = 22.59821958381127
$ws.Range("D24").Value = 4.747015846311729
$ws.Range("E24").Value = 35.13865545386768
$ws.Range("F24").Value = 20.51479855174997
$ws.Range("H24").Value = 4.049291792696274
$ws.Range("I24").Value = 3.857953781632018
$ws.Range("J24").ClearContents()
$ws.Range("O24").Value = 0
$ws.Range("P24").Value = 0
$ws.Range("Q24").Value = 15.29740027146494

# Row 25
$ws.Range("B25").Value = 26.72027819732691
$ws.Range("C25").Value = 20.69398490198719
$ws.Range("D25").Value = 4.500641749950556
$ws.Range("E25").Value = 31.88598415281427
$ws.Range("F25").Value = 19.25577469052761
$ws.Range("H25").Value = 3.57136126208584
$ws.Range("I25").Value = 3.489440049676813
$ws.Range("J25").ClearContents()
$ws.Range("O25").Value = 0
$ws.Range("P25").Value = 0
$ws.Range("Q25").Value = 14.71302920921626

